$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 79, pushing the existing rows 79-103 down to 80-104.
$ws.Rows("79:79").Insert()

# Populate the newly inserted row 79 with the new record.
$ws.Cells.Item(79, 1).Value = 1
$ws.Cells.Item(79, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(79, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(79, 4).Value = 44736
$ws.Cells.Item(79, 5).Value = 15
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100102
$ws.Cells.Item(79, 8).Value = "Cítricos"
$ws.Cells.Item(79, 9).Value = 100102004
$ws.Cells.Item(79, 10).Value = "Mandarina"
$ws.Cells.Item(79, 11).Value = "Clemenuless"
$ws.Cells.Item(79, 12).Value = "Tercera"
$ws.Cells.Item(79, 13).Value = 270
$ws.Cells.Item(79, 14).Value = 10000
$ws.Cells.Item(79, 15).Value = 11000
$ws.Cells.Item(79, 16).Value = 10500
$ws.Cells.Item(79, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(79, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(79, 19).Value = 525
$ws.Cells.Item(79, 20).Value = 20
